$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 14:21"

# Row 6
$ws.Range("B6").Value = 2719499
$ws.Range("C6").Value = 17895
$ws.Range("D6").Value = 1992150
$ws.Range("E6").Value = 675289
$ws.Range("G6").Value = 135
$ws.Range("H6").Value = 52060

# Row 41
$ws.Range("B41").Value = 77470
$ws.Range("C41").Value = 643
$ws.Range("D41").Value = 69243
$ws.Range("E41").Value = 7722
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 505

# Row 76
$ws.Range("A76").Value = "Estado de Palestina"
$ws.Range("B76").Value = 17306
$ws.Range("C76").Value = 462
$ws.Range("D76").Value = 9939
$ws.Range("E76").Value = 7254
$ws.Range("H76").Value = 113

# Row 77
$ws.Range("A77").Value = "Costa de Marfil"
$ws.Range("B77").Value = 17107
$ws.Range("D77").Value = 13990
$ws.Range("E77").Value = 3007
$ws.Range("H77").Value = 110

# Row 78
$ws.Range("B78").Value = 16351
$ws.Range("C78").Value = 240
$ws.Range("D78").Value = 10279
$ws.Range("E78").Value = 5577
$ws.Range("G78").Value = 9
$ws.Range("H78").Value = 495

# Row 79
$ws.Range("A79").Value = "Dinamarca"
$ws.Range("B79").Value = 15855
$ws.Range("C79").Value = 115
$ws.Range("D79").Value = 13529
$ws.Range("E79").Value = 1705
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 621

# Row 80
$ws.Range("A80").Value = "Corea del Sur"
$ws.Range("B80").Value = 15761
$ws.Range("C80").Value = 246
$ws.Range("D80").Value = 13934
$ws.Range("E80").Value = 1521
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 306

# Row 88
$ws.Range("B88").Value = 9981
$ws.Range("C88").Value = 142
$ws.Range("D88").Value = 8776
$ws.Range("E88").Value = 941

# Row 89
$ws.Range("B89").Value = 9721
$ws.Range("C89").Value = 15
$ws.Range("D89").Value = 8882
$ws.Range("E89").Value = 596

# Row 102
$ws.Range("A102").Value = "Croacia"
$ws.Range("B102").Value = 6855
$ws.Range("C102").Value = 199
$ws.Range("D102").Value = 5318
$ws.Range("E102").Value = 1371
$ws.Range("H102").Value = 166

# Row 103
$ws.Range("A103").Value = "Mauritania"
$ws.Range("B103").Value = 6762
$ws.Range("D103").Value = 6018
$ws.Range("E103").Value = 587
$ws.Range("H103").Value = 157

# Row 125
$ws.Range("D125").Value = 2755
$ws.Range("E125").Value = 134

# Row 158
$ws.Range("B158").Value = 989
$ws.Range("C158").Value = 6
$ws.Range("D158").Value = 520
$ws.Range("E158").Value = 444

# Row 165
$ws.Range("B165").Value = 702
$ws.Range("C165").Value = 3
$ws.Range("E165").Value = 3

# Row 213
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# Row 216
$ws.Range("A216").Value = "Islas Virgenes Britanicas"
$ws.Range("B216").Value = 11
$ws.Range("C216").Value = 2
$ws.Range("E216").Value = 2

# Row 217
$ws.Range("A217").Value = "Sahara Occidental"
$ws.Range("B217").Value = 10
$ws.Range("D217").Value = 8
